# Tata Motors demerger: split the single TATAMOTORS holding into two
# successor positions - TMCV (Tata Motors Limited) and TMPV (Tata Motors
# Passenger Vehicle Limited) - each carried at its own purchase price while
# keeping every other holding intact (just shifted down one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-mark the existing Purchase_Date cell as Text *before* the row insert so
# that when we retype it below it is stored as the literal string
# "2024-01-15" (reusing the existing shared string) instead of being
# auto-parsed into a date serial.
$ws.Range("C2").NumberFormatLocal = "@"

# Make room for the new TMPV row: push the old row 3 (and everything below
# it) down by one.
$ws.Rows("3:3").Insert()

# Row 2 becomes TMCV - Tata Motors Limited. Purchase_Date stays the plain
# text "2024-01-15" and Total_Investment is a hard-coded number (not a
# formula), matching how the split was recorded.
$ws.Range("A2").Value = "TMCV"
$ws.Range("B2").Value = "Tata Motors Limited"
$ws.Range("C2").Value = "2024-01-15"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = 124.11
$ws.Range("E2").Value = 12
$ws.Range("F2").Value = 1489.31

# New row 3 is TMPV - Tata Motors Passenger Vehicle Limited. This one's
# Purchase_Date is entered as an actual date value with a custom date
# number format, and Total_Investment is likewise a hard-coded number.
$ws.Range("A3").Value = "TMPV"
$ws.Range("B3").Value = "Tata Motors Pass Vehicle Limited"
$ws.Range("C3").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("C3").Value = 45306
$ws.Range("D3").Value = 274.32
$ws.Range("E3").Value = 12
$ws.Range("F3").Value = 3291.78

# Restore the Total_Investment formula (Purchase_Price * Quantity) for every
# surviving holding, now living one row further down (rows 4-16).
for ($r = 4; $r -le 16; $r++) {
    $ws.Range("F$r").Formula = "=D$r*E$r"
}

# Purchase_Date column widened to fit the dates again.
$ws.Columns("C").ColumnWidth = 9.33

# Leave the selection where the author ended up.
$ws.Range("F19").Select()
